$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 41, shifting existing rows 41-43 down to 42-44
$ws.Rows.Item(41).Insert()

# Copy style from the row above (row 40) into the new row 41 so formatting matches
$ws.Range("A40:D40").Copy()
$ws.Range("A41:D41").PasteSpecial(-4122) # xlPasteFormats

# Populate the new row 41 with the inserted record
$ws.Range("A41").Value = "Z16_B04"
$ws.Range("B41").Value = "Z16"
$ws.Range("C41").Value = "Finanzkriminalität"
$ws.Range("D41").Value = "XXXFinanzkriminalität"
